# Fruta / hortaliza, semanal
# Rotate the weekly price-record rows: row2 <- old row4, row3 <- old row2, row4 <- old row3
# (only the columns that actually differ between the three source rows need touching;
#  columns that are identical across rows are left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that vary row to row.
# (.Value2 is used for reads -- .Value's getter is unreliable in this host)
$D2 = $ws.Range("D2").Value2
$M2 = $ws.Range("M2").Value2
$N2 = $ws.Range("N2").Value2
$O2 = $ws.Range("O2").Value2
$P2 = $ws.Range("P2").Value2
$R2 = $ws.Range("R2").Value2
$S2 = $ws.Range("S2").Value2

$D3 = $ws.Range("D3").Value2
$M3 = $ws.Range("M3").Value2
$N3 = $ws.Range("N3").Value2
$O3 = $ws.Range("O3").Value2
$P3 = $ws.Range("P3").Value2
$R3 = $ws.Range("R3").Value2
$S3 = $ws.Range("S3").Value2

$D4 = $ws.Range("D4").Value2
$M4 = $ws.Range("M4").Value2
$N4 = $ws.Range("N4").Value2
$O4 = $ws.Range("O4").Value2
$P4 = $ws.Range("P4").Value2
$R4 = $ws.Range("R4").Value2
$S4 = $ws.Range("S4").Value2

# Row 2 becomes old row 4
$ws.Range("D2").Value2 = $D4
$ws.Range("M2").Value2 = $M4
$ws.Range("N2").Value2 = $N4
$ws.Range("O2").Value2 = $O4
$ws.Range("P2").Value2 = $P4
$ws.Range("R2").Value2 = $R4
$ws.Range("S2").Value2 = $S4

# Row 3 becomes old row 2
$ws.Range("D3").Value2 = $D2
$ws.Range("M3").Value2 = $M2
$ws.Range("N3").Value2 = $N2
$ws.Range("O3").Value2 = $O2
$ws.Range("P3").Value2 = $P2
$ws.Range("R3").Value2 = $R2
$ws.Range("S3").Value2 = $S2

# Row 4 becomes old row 3
$ws.Range("D4").Value2 = $D3
$ws.Range("M4").Value2 = $M3
$ws.Range("N4").Value2 = $N3
$ws.Range("O4").Value2 = $O3
$ws.Range("P4").Value2 = $P3
$ws.Range("R4").Value2 = $R3
$ws.Range("S4").Value2 = $S3
